$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Update user cell (A2) and email cell (B2) on Planilha1
$ws1.Range("A2").Value = "RodrigoVil55"
$ws1.Range("B2").Value = "rodrigovilnateste@rstinet.com"
$null = $ws1.Hyperlinks.Add($ws1.Range("B2"), "mailto:rodrigovilnateste@rstinet.com")

# Clear the leftover "screenshot" rows on Planilha1 (kept on Planilha2 only)
$ws1.Range("A6").ClearContents()
$ws1.Range("A7").ClearContents()
$ws1.Rows("7").AutoFit()

# Selection / active tab changes
$null = $ws1.Select()
$null = $ws1.Range("B5").Select()
